$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 32, shifting existing rows 32-35 down to 33-36.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly price record.
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = 44753
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 100112013
$ws.Range("G32").Value = "Alcachofa"
$ws.Range("H32").Value = "Argentina(o)"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 60
$ws.Range("K32").Value = 16000
$ws.Range("L32").Value = 17000
$ws.Range("M32").Value = 16500
$ws.Range("N32").Value = "`$/caja 50 unidades"
$ws.Range("O32").Value = "Provincia de Limarí"
$ws.Range("P32").Value = 330
$ws.Range("Q32").Value = 50
$ws.Range("R32").Value = "Hortaliza"
